$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("numeric")
$ws.Activate()

# Set E2 to the string value "*"
$ws.Range("E2").Value = "*"

# Update the selection to E2 (matches sqref/activeCell in diff)
$ws.Range("E2").Select()
